$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Справочники")

# Country reference list (column H) after removing "СЕВЕРНАЯ КОРЕЯ" and
# renaming a few entries; the list shifts up by one row and the final
# row (H243) becomes empty.
$countries = @(
    "АБХАЗИЯ",
    "АВСТРАЛИЯ",
    "АВСТРИЯ",
    "АЗЕРБАЙДЖАН",
    "АЛАНДСКИЕ ОСТРОВА",
    "АЛБАНИЯ",
    "АЛЖИР",
    "АМЕРИКАНСКИЕ ВИРГИНСКИЕ ОСТРОВА",
    "АМЕРИКАНСКОЕ САМОА",
    "АНГИЛЬЯ",
    "АНГОЛА",
    "АНДОРРА",
    "АНТАРКТИДА",
    "АНТИГУА И БАРБУДА",
    "АРГЕНТИНА",
    "АРМЕНИЯ",
    "АРУБА",
    "АФГАНИСТАН",
    "БАГАМЫ",
    "БАНГЛАДЕШ",
    "БАРБАДОС",
    "БАХРЕЙН",
    "БЕЛИЗ",
    "БЕЛАРУСЬ",
    "БЕЛЬГИЯ",
    "БЕНИН",
    "БЕРМУДЫ",
    "БОЛГАРИЯ",
    "БОЛИВИЯ",
    "БОСНИЯ И ГЕРЦЕГОВИНА",
    "БОТСВАНА",
    "БРАЗИЛИЯ",
    "БРИТАНСКАЯ ТЕРРИТОРИЯ В ИНДИЙСКОМ ОКЕАНЕ",
    "БРИТАНСКИЕ ВИРГИНСКИЕ ОСТРОВА",
    "БРУНЕЙ ДАРУССАЛАМ",
    "БУРКИНА-ФАСО",
    "БУРУНДИ",
    "БУТАН",
    "ВАНУАТУ",
    "ВАТИКАНСКОЕ ГОСУДАРСТВО (ПАПСКИЙ ПРЕСТОЛ)",
    "СОЕДИНЕННОЕ КОРОЛЕВСТВО",
    "ВЕНГРИЯ",
    "ВЕНЕСУЭЛА",
    "ВНЕШНИЕ МАЛЫЕ ОСТРОВА (США)",
    "ВОСТОЧНЫЙ ТИМОР",
    "ВЬЕТНАМ",
    "ГАБОН",
    "ГАИТИ",
    "ГАМБИЯ",
    "ГАНА",
    "ГВАДЕЛУПА",
    "ГВАТЕМАЛА",
    "ГВИНЕЯ",
    "ГВИНЕЯ-БИСАУ",
    "ГЕРМАНИЯ",
    "ГИБРАЛТАР",
    "ГОНДУРАС",
    "ГОНКОНГ",
    "ГРЕНАДА",
    "ГРЕНЛАНДИЯ",
    "ГРЕЦИЯ",
    "ГРУЗИЯ",
    "ГУАМ",
    "ДАНИЯ",
    "ДЕМОКРАТИЧЕСКАЯ РЕСПУБЛИКА КОНГО",
    "ДЖИБУТИ",
    "ДОМИНИКАНА",
    "ДОМИНИКАНСКАЯ РЕСПУБЛИКА",
    "ЕГИПЕТ",
    "ЗАМБИЯ",
    "ЗАПАДНАЯ САХАРА",
    "ЗИМБАБВЕ",
    "ИЗРАИЛЬ",
    "ИНДИЯ",
    "ИНДОНЕЗИЯ",
    "ИОРДАНИЯ",
    "ИРАК",
    "ИРАН (ИСЛАМСКАЯ РЕСПУБЛИКА)",
    "ИРЛАНДИЯ",
    "ИСЛАНДИЯ",
    "ИСПАНИЯ",
    "ИТАЛИЯ",
    "ЙЕМЕН",
    "КАБО-ВЕРДЕ",
    "КАЗАХСТАН",
    "КАЙМАНОВЫ ОСТРОВА",
    "КАМБОДЖА",
    "КАМЕРУН",
    "КАНАДА",
    "КАТАР",
    "КЕНИЯ",
    "КИПР",
    "КИРГИЗИЯ",
    "КИРИБАТИ",
    "КИТАЙ",
    "КОКОСОВЫЕ ОСТРОВА (ОСТРОВА КИЛИНГ)",
    "КОЛУМБИЯ",
    "КОМОРСКИЕ ОСТРОВА",
    "КОСТА-РИКА",
    "КОТ-Д'ИВУАР",
    "КУБА",
    "КУВЕЙТ",
    "ЛАОССКАЯ НАРОДНАЯ ДЕМОКРАТИЧЕСКАЯ РЕСПУБЛИКА",
    "ЛАТВИЯ",
    "ЛЕСОТО",
    "ЛИБЕРИЯ",
    "ЛИВАН",
    "ЛИВИЯ",
    "ЛИТВА",
    "ЛИХТЕНШТЕЙН",
    "ЛЮКСЕМБУРГ",
    "МАВРИКИЙ",
    "МАВРИТАНИЯ",
    "МАДАГАСКАР",
    "МАЙОТТА",
    "МАКАО",
    "МАЛАВИ",
    "МАЛАЙЗИИЯ",
    "МАЛИ",
    "МАЛЬДИВЫ",
    "МАЛЬТА",
    "МАРОККО",
    "МАРТИНИКА",
    "МАРШАЛЛОВЫ ОСТРОВА",
    "МЕКСИКА",
    "МИКРОНЕЗИЯ, ОБЪЕДИНЕННЫЕ ГОСУДАРСТВА",
    "МОЗАМБИК",
    "МОЛДАВИЯ",
    "МОНАКО",
    "МОНГОЛИЯ",
    "МОНТСЕРРАТ",
    "МЬЯНМА",
    "НАМИБИЯ",
    "НАУРУ",
    "НЕПАЛ",
    "НИГЕР",
    "НИГЕРИЯ",
    "НИДЕРЛАНДЫ",
    "НИДЕРЛАНДЫ АНТИЛЬСКИЕ ОСТРОВА",
    "НИКАРАГУА",
    "НИУЭ",
    "НОВАЯ ЗЕЛАНДИЯ",
    "НОВАЯ КАЛЕДОНИЯ",
    "НОРВЕГИЯ",
    "ОБЪЕДИНЕННЫЕ АРАБСКИЕ ЭМИРАТЫ",
    "ОМАН",
    "ОСТРОВ БУВЕ",
    "ОСТРОВ НОРФОЛК",
    "ОСТРОВ РОЖДЕСТВА",
    "ОСТРОВ СВЯТОЙ ЕЛЕНЫ",
    "ОСТРОВА КУКА",
    "ОСТРОВА ПИТКЭРН",
    "ПАКИСТАН",
    "ПАЛАУ",
    "ПАЛЕСТИНА",
    "ПАНАМА",
    "ПАПУА-НОВАЯ ГВИНЕЯ",
    "ПАРАГВАЙ",
    "ПЕРУ",
    "ПОЛЬША",
    "ПОРТУГАЛИЯ",
    "ПУЭРТО-РИКО",
    "РЕСПУБЛИКА КОНГО",
    "РЕСПУБЛИКА КОРЕЯ",
    "РЕСПУБЛИКА МАКЕДОНИЯ, БЫВШАЯ ЮГОСЛАВСКАЯ РЕСПУБЛИКА",
    "РЕЮНЬОН",
    "РОССИЯ",
    "РУАНДА",
    "РУМЫНИЯ",
    "САЛЬВАДОР",
    "САМОА",
    "САН-МАРИНО",
    "САН-ТОМЕ И ПРИНСИПИ",
    "САУДОВСКАЯ АРАВИЯ",
    "СВАЗИЛЕНД",
    "СЕВЕРНЫЕ МАРИАНСКИЕ ОСТРОВА",
    "СЕЙШЕЛЬСКИЕ ОСТОРВА",
    "СЕН-ПЬЕР И МИКЕЛОН",
    "СЕНЕГАЛ",
    "СЕНТ-ВИНСЕНТ И ГРЕНАДИНЫ",
    "СЕНТ-КИТТС И НЕВИС",
    "СЕНТ-ЛЮСИЯ",
    "СЕРБИЯ",
    "СИНГАПУР",
    "СИРИЯ",
    "СЛОВАКИЯ",
    "СЛОВЕНИЯ",
    "СОЕДИНЕННЫЕ ШТАТЫ",
    "СОЛОМОНОВЫ ОСТРОВА",
    "СОМАЛИ",
    "СУДАН",
    "СУРИНАМ",
    "СЬЕРРА-ЛЕОНЕ",
    "ТАДЖИКИСТАН",
    "ТАИЛАНД",
    "ТАЙВАНЬ, ПРОВИНЦИЯ КИТАЯ",
    "ТАНЗАНИЯ, ОБЪЕДИНЕННАЯ РЕСПУБЛИКА",
    "ТЕРКС И КАЙКОС",
    "ТОГО",
    "ТОКЕЛАУ",
    "ТОНГА",
    "ТРИНИДАД И ТОБАГО",
    "ТУВАЛУ",
    "ТУНИС",
    "ТУРКМЕНИЯ",
    "ТУРЦИЯ",
    "УГАНДА",
    "УЗБЕКИСТАН",
    "УКРАИНА",
    "УОЛЛИС И ФУТУНА",
    "УРУГВАЙ",
    "ФАРЕРСКИЕ ОСТРОВА",
    "ФИДЖИ",
    "ФИЛИППИНЫ",
    "ФИНЛЯНДИЯ",
    "ФОЛКЛЕНДСКИЕ ОСТРОВА (МАЛЬВИНСКИЕ)",
    "ФРАНЦИЯ",
    "ФРАНЦУЗСКАЯ ГВИАНА",
    "ФРАНЦУЗСКАЯ ПОЛИНЕЗИЯ",
    "ФРАНЦУЗСКИЕ ЮЖНЫЕ ТЕРРИТОРИИ",
    "ХЕРДА И МАКДОНАЛЬДА ОСТРОВА",
    "ХОРВАТИЯ",
    "ЦЕНТРАЛЬНОАФРИКАНСКАЯ РЕСПУБЛИКА",
    "ЧАД",
    "ЧЕРНОГОРИЯ",
    "ЧЕХИЯ",
    "ЧИЛИ",
    "ШВЕЙЦАРИЯ",
    "ШВЕЦИЯ",
    "ШПИЦБЕРГЕН И ЯН-МАЙЕН",
    "ШРИ ЛАНКА",
    "ЭКВАДОР",
    "ЭКВАТОРИАЛЬНАЯ ГВИНЕЯ",
    "ЭРИТРЕЯ",
    "ЭСТОНИЯ",
    "ЭФИОПИЯ",
    "ЮЖНАЯ АФРИКА",
    "ЮЖНАЯ ГЕОРГИЯ И ЮЖНЫЕ САНДВИЧЕВЫ ОСТРОВА",
    "ЮЖНАЯ ОСЕТИЯ",
    "ЯМАЙКА",
    "ЯПОНИЯ"
)

for ($i = 0; $i -lt $countries.Length; $i++) {
    $ws.Cells.Item($i + 2, 8).Value = $countries[$i]
}

$ws.Cells.Item(243, 8).ClearContents()
